# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45941

$ws.Range("B2").Value = 106.95
$ws.Range("C2").Value = 105.45
$ws.Range("D2").Value = 103.81
$ws.Range("E2").Value = 101.03
$ws.Range("F2").Value = 100.34
$ws.Range("G2").Value = 101.89
$ws.Range("H2").Value = 102.35
$ws.Range("I2").Value = 101.41
$ws.Range("J2").Value = 102.09
$ws.Range("K2").Value = 85.69
$ws.Range("L2").Value = 53.84
$ws.Range("M2").Value = 18.25
$ws.Range("N2").Value = 4.88
$ws.Range("O2").Value = 4.44
$ws.Range("P2").Value = 3.84
$ws.Range("Q2").Value = 3.72
$ws.Range("R2").Value = 7.51
$ws.Range("S2").Value = 41.98
$ws.Range("T2").Value = 83.06
$ws.Range("U2").Value = 106.6
$ws.Range("V2").Value = 108.26
$ws.Range("W2").Value = 109.81
$ws.Range("X2").Value = 106.99
$ws.Range("Y2").Value = 104.65
$ws.Range("Z2").Value = 73.7

$ws.Range("AB2").Value = 107.43
$ws.Range("AD2").Value = 109.04
$ws.Range("AE2").Value = "0h-2h"
$ws.Range("AF2").Value = 106.2
$ws.Range("AG2").Value = "10h-17h"
